$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correccion de matriz de casos de prueba
# Remove the word "izquierdo" (and tidy spacing) from the "Dar click" instructions
# in column D for rows 3-8. The assignment order below matches the order in
# which the author's edit re-introduced these strings into the shared string
# table (Zologico Zacano, Lugares Turisticos, mouse/Volcan, Eventos, Historia,
# Bolsa de trabajo).

$ws.Range("D4").Value = "1.- Ingresar a la página incial`n2.- Posicionar el mause sobre la palabra ""Lugares Turisticos""`n3.- En linea recta verticalmente, posicionar el mause en el nombre ""Zologico Zacano""`n4.- Dar click en ""Zologico Zacano"""

$ws.Range("D3").Value = "1.- Ingresar a la página incial`n2.- Posicionar el mause sobre la palabra ""Lugares Turisticos""`n3.- Dar click  en ""Lugares Turisticos"""

$ws.Range("D5").Value = "1.- Posicionar el mouse sobre Lugares Turisticos.`n2.- Dar click izquierdo `n2.-Posicionar el mouse sobre Volcan popocatepetl.`n2.- Dar click  ."

$ws.Range("D6").Value = "1.- Ingresar a la página incial`n2.- Posicionar el mause sobre la palabra ""Eventos""`n3.- Dar click en ""Eventos"""

$ws.Range("D7").Value = "1.- Ingresar a la página incial`n2.- Posicionar el mause sobre la palabra ""Historia""`n3.- Dar click en ""Historia"""

$ws.Range("D8").Value = "1.- Ingresar a la página incial`n2.- Posicionar el mause sobre la palabra ""Bolsa de Trabajo""`n3.- Dar click en ""Bolsa de trabajo"""

# Update the saved view state (scroll position / active selection)
$ws.Range("D19").Select()
$excel.ActiveWindow.ScrollRow = 17
